$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "slide.types": insert a new "slide.layout" column after
# "slide.type.position" (old column C), before the old "slide.loop.var"
# column (old column D). This shifts slide.loop.var -> E and
# slide.graph.type -> F.
# ---------------------------------------------------------------------------
$wsTypes = $wb.Worksheets.Item("slide.types")
$wsTypes.Columns.Item(4).Insert()

$wsTypes.Range("D1").Value = "slide.layout"
$wsTypes.Range("D2").Value = "title"
$wsTypes.Range("D3").Value = "slide"
$wsTypes.Range("D4").Value = "slide"
$wsTypes.Range("D5").Value = "slide"
$wsTypes.Range("D6").Value = "section"
$wsTypes.Range("D7").Value = "slide"
$wsTypes.Range("D8").Value = "slide"
$wsTypes.Range("D9").Value = "section"
$wsTypes.Range("D10").Value = "slide"
$wsTypes.Range("D11").Value = "slide"
$wsTypes.Range("D12").Value = "section"
$wsTypes.Range("D13").Value = "slide"

$wsTypes.Columns.Item(4).ColumnWidth = $wsTypes.Columns.Item(3).ColumnWidth()

# ---------------------------------------------------------------------------
# Sheet "slide.objects": rename the "id" header (column B) to "object.id".
# ---------------------------------------------------------------------------
$wsObjects = $wb.Worksheets.Item("slide.objects")
$wsObjects.Range("B1").Value = "object.id"

# ---------------------------------------------------------------------------
# Selection / active-sheet bookkeeping to match the saved view state: the
# user had D14 selected on slide.types, then switched to slide.objects and
# selected B1 there (making it the active tab).
# ---------------------------------------------------------------------------
$wsTypes.Range("D14").Select()
$wsObjects.Range("B1").Select()
